$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = '66.248.35'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '3.541.28'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.65'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.19'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '3.539.92'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '8.09'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.136'
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = '  -4.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.410'
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('D13').Value = '4.151.02'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000207'
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = '  -4.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.08'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  -4.81%  '
$ws.Range('D16').Value = '3.545.72'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = '66.345.46'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.91'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.89'
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '425.66'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.601'
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.02'
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').Value = '3.686.32'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.09'
$ws.Range('D28').Style = $plainStyle
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.18'
$ws.Range('D29').Style = $plainStyle
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.49'
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.48'
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = '  -6.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.159'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.30'
$ws.Range('D34').Style = $plainStyle
$ws.Range('D35').Value = '3.535.02'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.75'
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  -3.13%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.81'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.62'
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '173.46'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0856'
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = '  -4.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.26'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.892'
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.90'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  -6.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.59'
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('E48').Value = '  -6.81%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.941'
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = '  -5.57%  '
